$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '62.128.88'
Set-TextValue "E2" '  -1.05%  '

Set-TextValue "D3" '3.422.97'
Set-TextValue "E3" '  -1.66%  '

Set-TextValue "E4" '  -0.13%  '

Set-TextValue "D5" '408.09'
Set-TextValue "E5" '  -2.28%  '

Set-TextValue "D6" '133.12'
Set-TextValue "E6" '  +1.51%  '

Set-TextValue "D7" '0.593'
Set-TextValue "E7" '  -0.56%  '

Set-TextValue "E8" '  -0.06%  '

Set-TextValue "D9" '0.676'
Set-TextValue "E9" '  -2.53%  '

Set-TextValue "E10" '  -3.44%  '

Set-TextValue "D11" '42.49'
Set-TextValue "E11" '  -3.56%  '

Set-TextValue "E12" '  -1.59%  '

Set-TextValue "D13" '3.963.27'
Set-TextValue "E13" '  -1.66%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D14" '19.98'
Set-TextValue "E14" '  -2.06%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D15" '8.45'
Set-TextValue "E15" '  -4.24%  '

Set-TextValue "D16" '3.422.21'
Set-TextValue "E16" '  -1.89%  '

Set-TextValue "D17" '62.086.59'
Set-TextValue "E17" '  -1.18%  '

Set-TextValue "E18" '  -3.46%  '

Set-TextValue "D19" '11.04'
Set-TextValue "E19" '  +0.17%  '

Set-TextValue "E20" '  -3.23%  '

Set-TextValue "E21" '  -4.37%  '

Set-TextValue "D22" '85.02'
Set-TextValue "E22" '  +3.45%  '

Set-TextValue "D23" '316.20'
Set-TextValue "E23" '  +0.02%  '

Set-TextValue "E24" '  -3.37%  '

Set-TextValue "E25" '  -3.58%  '

Set-TextValue "E26" '  +9.19%  '

Set-TextValue "D27" '29.78'
Set-TextValue "E27" '  -3.47%  '

Set-TextValue "D28" '8.23'
Set-TextValue "E28" '  +1.28%  '

Set-TextValue "D29" '7.69'
Set-TextValue "E29" '  -1.12%  '

Set-TextValue "D30" '2.74'
Set-TextValue "E30" '  +3.30%  '

Set-TextValue "D31" '0.174'
Set-TextValue "E31" '  -2.56%  '

Set-TextValue "E32" '  -5.16%  '

Set-TextValue "D33" '43.03'
Set-TextValue "E33" '  -3.56%  '

Set-TextValue "D34" '11.43'
Set-TextValue "E34" '  -4.52%  '

Set-TextValue "E35" '  -0.09%  '

Set-TextValue "D36" '0.0487'
Set-TextValue "E36" '  -2.34%  '

Set-TextValue "D37" '52.11'
Set-TextValue "E37" '  -0.81%  '

Set-TextValue "E38" '  +0.19%  '

Set-TextValue "E39" '  -4.65%  '

Set-TextValue "D40" '3.00'
Set-TextValue "E40" '  -1.04%  '

Set-TextValue "E41" '  -0.96%  '

Set-TextValue "D42" '138.18'
Set-TextValue "E42" '  +0.50%  '

Set-TextValue "D43" '0.125'
Set-TextValue "E43" '  -1.09%  '

Set-TextValue "E44" '  +1.25%  '

Set-TextValue "D45" '4.00'
Set-TextValue "E45" '  -1.31%  '

Set-TextValue "D46" '16.87'
Set-TextValue "E46" '  -4.25%  '

Set-TextValue "E47" '  -3.06%  '

Set-TextValue "D48" '21.51'
Set-TextValue "E48" '  -6.15%  '

Set-TextValue "D49" '2.135.07'
Set-TextValue "E49" '  -5.39%  '

Set-TextValue "E50" '  -4.48%  '

Set-TextValue "D51" '1.89'
Set-TextValue "E51" '  +0.20%  '
